$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "RM 232" row (row 26) and the "SC 92" row (now row 27 after the
# first deletion) from the dataset. Every row below shifts up accordingly.
$ws.Rows(26).Delete()
$ws.Rows(27).Delete()

# Re-mask / fill values for the "Error Calculations and Plots" pass over
# column F (and a couple of others), matching the new row positions.
$ws.Range("D3").Value = ""
$ws.Range("F5").Value = ""
$ws.Range("F8").Value = 17.05
$ws.Range("F10").Value = 16.43
$ws.Range("F12").Value = ""
$ws.Range("F15").Value = 16.2
$ws.Range("F18").Value = ""
$ws.Range("F19").Value = ""
$ws.Range("F25").Value = 16.6
$ws.Range("B26").Value = -20.2
$ws.Range("B27").Value = ""
$ws.Range("F29").Value = ""
$ws.Range("B33").Value = -19.5
$ws.Range("D33").Value = -14.1
